$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Pairwise row swaps (columns B:AC; column A index stays put) ---
$swapPairs = @(
    @(41,42),
    @(77,78),
    @(79,80),
    @(117,119),
    @(123,124),
    @(126,127),
    @(152,153),
    @(159,160),
    @(169,171),
    @(208,209),
    @(211,212)
)

foreach ($pair in $swapPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    $v1 = $ws.Range("B$($r1):AC$($r1)").Value()
    $v2 = $ws.Range("B$($r2):AC$($r2)").Value()
    $ws.Range("B$($r1):AC$($r1)").Value = $v2
    $ws.Range("B$($r2):AC$($r2)").Value = $v1
}

# --- 3-way rotation: new47 = old49, new48 = old47, new49 = old48 ---
$v47 = $ws.Range("B47:AC47").Value()
$v48 = $ws.Range("B48:AC48").Value()
$v49 = $ws.Range("B49:AC49").Value()
$ws.Range("B47:AC47").Value = $v49
$ws.Range("B48:AC48").Value = $v47
$ws.Range("B49:AC49").Value = $v48

# --- Shift rows 215..219 up from 217..221, then delete old trailing rows 220,221 ---
$sv215 = $ws.Range("B217:AC217").Value()
$sv216 = $ws.Range("B218:AC218").Value()
$sv217 = $ws.Range("B219:AC219").Value()
$sv218 = $ws.Range("B220:AC220").Value()
$sv219 = $ws.Range("B221:AC221").Value()

$ws.Range("B215:AC215").Value = $sv215
$ws.Range("B216:AC216").Value = $sv216
$ws.Range("B217:AC217").Value = $sv217
$ws.Range("B218:AC218").Value = $sv218
$ws.Range("B219:AC219").Value = $sv219

$ws.Rows("220:221").Delete()
